$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (stiction torque) re-recorded for rows 2-11
$ws.Range("C2").Value = 6.361236
$ws.Range("C3").Value = 6.361236
$ws.Range("C4").Value = 6.361236
$ws.Range("C5").Value = 6.361236
$ws.Range("C6").Value = 6.361236
$ws.Range("C7").Value = 6.361236
$ws.Range("C8").Value = 6.557526
$ws.Range("C9").Value = 6.655671
$ws.Range("C10").Value = 6.901033
$ws.Range("C11").Value = 7.048250

# Column B updated alongside the re-recording (row 2 unchanged)
$ws.Range("B3").Value = -0.007835
$ws.Range("B4").Value = -0.012120
$ws.Range("B5").Value = -0.016405
$ws.Range("B6").Value = -0.020690
$ws.Range("B7").Value = -0.024975
$ws.Range("B8").Value = -0.029260
$ws.Range("B9").Value = -0.033545
$ws.Range("B10").Value = -0.037830
$ws.Range("B11").Value = -0.042115

# Column D (report part) recalculated values
$ws.Range("D8").Value = 0.196290
$ws.Range("D9").Value = 0.294435
$ws.Range("D10").Value = 0.539797
$ws.Range("D11").Value = 0.687015
